$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions): reset "想去人数" (want-to-go count) column F, rows 2-23
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2:F23").Value = 0

# Sheet "演出" (Performances): reset column F, row 2
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2:F2").Value = 0

# Sheet "全部类型" (All types): reset column F, rows 2-24
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2:F24").Value = 0
